$wb = $excel.ActiveWorkbook

# Both "展览" and "全部类型" sheets carry duplicated event listings;
# refresh the "想去人数" (interested-count) figures for the two events
# that changed in this scrape run.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 499
    $ws.Range("F6").Value = 685
}
